$d = $word.ActiveDocument

# Turn off "smart quotes" / autoformat-as-you-type so Find & Replace does not
# mangle straight/curly quote characters that happen to sit near our matches.
$word.Options.AutoFormatAsYouTypeReplaceQuotes = $false
$word.Options.AutoFormatReplaceQuotes = $false

# wdReplaceAll = 2 ; wdFindContinue = 1
$wdReplaceAll = 2
$wdFindContinue = 1

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $replace, $wdReplaceAll) | Out-Null
}

# 1. Wrap {{ plaintiffs }} with capitalize(...)
Replace-Text "{{ plaintiffs }}" "{{ capitalize(plaintiffs) }}"

# 2. Wrap {{ defendants }} with capitalize(...)
Replace-Text "{{ defendants }}" "{{ capitalize(defendants) }}"

# 3. Remove the stray space between "%}" and "{{" in the judgment_appealing if block
#    (kept narrow so the surrounding 'Other' quotes are left completely untouched)
Replace-Text "%} {{ judgment_appealing" "%}{{ judgment_appealing"

# 4. Remove the non-breaking space between "{% endif %}" and the following "{{ showifdef("
$nbsp = [char]0x00A0
$findNbsp = "endif %}" + $nbsp + "{{ "
Replace-Text $findNbsp "endif %}{{ "

# 5. Rename judgment_date -> trial_court_judgment_date
Replace-Text "on {{ judgment_date }}." "on {{ trial_court_judgment_date }}."

# 6. "him/her" -> "them"
Replace-Text "him/her" "them"

# 7. "his/her lawyer" -> "their lawyer"
Replace-Text "his/her lawyer" "their lawyer"
